$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "Haunted Carriage 12PM"
$ws.Range("A6").Value = "Haunted Carriage 830PM"
$ws.Range("A7").Value = "Haunted Carriage 10PM"

$ws.Range("C8").Select()
